# Apply the update described by the diff:
#  - Row 10: several field updates (id, validation status, locality name,
#    accuracy, start/end date, public comment, "not refound" flag,
#    reporter/observers)
#  - New rows 16, 17, 18 appended with fresh observation records

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force a text-looking value (e.g. a date string like "2014-04-15")
    # to be stored as literal text instead of being auto-parsed into a
    # serial date number, matching how the source data is encoded. A
    # leading apostrophe is Excel's standard "treat as text" input prefix.
    $ws.Cells.Item($row, $col).Value2 = "'" + $text
}

function Set-EmptyTextCell($row, $col) {
    # Create a present-but-blank, text-typed cell (mirrors the source
    # file's empty <is><t/></is> placeholder cells) without altering the
    # used range. A lone apostrophe is stored as an empty text value.
    $ws.Cells.Item($row, $col).Value2 = "'"
}

# ---------------------------------------------------------------------
# Row 10 updates
# ---------------------------------------------------------------------
$ws.Cells.Item(10, 1).Value2 = 82402784                                   # A10 Id
$ws.Cells.Item(10, 3).Value2 = "Behöver inte valideras"                   # C10 Valideringsstatus
$ws.Cells.Item(10, 16).Value2 = "59, Öl"                                  # P10 Lokalnamn
$ws.Cells.Item(10, 19).Value2 = 50                                        # S10 Noggrannhet
Set-TextCell 10 25 "2014-04-15"                                           # Y10 Startdatum
Set-TextCell 10 27 "2014-04-15"                                           # AA10 Slutdatum
$ws.Cells.Item(10, 29).Value2 = "Uttorkat, inget vatten"                  # AC10 Publik kommentar
$ws.Cells.Item(10, 30).Value2 = $true                                     # AD10 Ej återfunnen
$ws.Cells.Item(10, 49).Value2 = "Pia Hertonsson"                          # AW10 Rapportör
$ws.Cells.Item(10, 50).Value2 = "Pia Hertonsson, Marika Stenberg, Erik Fridolf, Per Nyström, Lars-Göran Pärlklint"  # AX10 Observatörer

# ---------------------------------------------------------------------
# New rows 16-18
# ---------------------------------------------------------------------
$newRows = @(
    @{
        Row=16; A=112272267; B=98961; C="Ovaliderad"; D="LC"; E=222498
        F="Blåsippa"; G="Hepatica nobilis"; H="Schreb."
        P="Störlinge gran, Öl"; Q=607474; R=6297135; S=25
        T="Kalmar"; U="Borgholm"; V="Öland"; W="Gärdslösa"
        Y="2023-09-23"; AA="2023-09-23"
        AW="Ulla-Britt Andersson"; AX="Ulla-Britt Andersson, Thomas Gunnarsson"
    },
    @{
        Row=17; A=112272269; B=96698; C="Ovaliderad"; D="LC"; E=219798
        F="Skogsknipprot"; G="Epipactis helleborine"; H="(L.) Crantz"
        P="Störlinge gran, Öl"; Q=607474; R=6297135; S=25
        T="Kalmar"; U="Borgholm"; V="Öland"; W="Gärdslösa"
        Y="2023-09-23"; AA="2023-09-23"
        AW="Ulla-Britt Andersson"; AX="Ulla-Britt Andersson, Thomas Gunnarsson"
    },
    @{
        Row=18; A=112272399; B=108537; C="Ovaliderad"; D="LC"; E=219677
        F="Murgröna"; G="Hedera helix"; H="L."
        P="Störlinge gran, Öl"; Q=607628; R=6297526; S=25
        T="Kalmar"; U="Borgholm"; V="Öland"; W="Gärdslösa"
        Y="2023-09-23"; AA="2023-09-23"
        AW="Ulla-Britt Andersson"; AX="Ulla-Britt Andersson, Thomas Gunnarsson"
    }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.A        # Id
    $ws.Cells.Item($row, 2).Value2 = $r.B        # Taxonsorteringsordning
    $ws.Cells.Item($row, 3).Value2 = $r.C        # Valideringsstatus
    $ws.Cells.Item($row, 4).Value2 = $r.D        # Rödlistade
    $ws.Cells.Item($row, 5).Value2 = $r.E        # TaxonId
    $ws.Cells.Item($row, 6).Value2 = $r.F        # Artnamn
    $ws.Cells.Item($row, 7).Value2 = $r.G        # Vetenskapligt namn
    $ws.Cells.Item($row, 8).Value2 = $r.H        # Auktor
    Set-EmptyTextCell $row 9                     # Antal (blank placeholder)
    Set-EmptyTextCell $row 10                    # Enhet (blank placeholder)
    Set-EmptyTextCell $row 11                    # Ålder-Stadium (blank placeholder)
    Set-EmptyTextCell $row 12                    # Kön (blank placeholder)
    Set-EmptyTextCell $row 14                    # Metod (blank placeholder)
    $ws.Cells.Item($row, 16).Value2 = $r.P       # Lokalnamn
    $ws.Cells.Item($row, 17).Value2 = $r.Q       # Ost
    $ws.Cells.Item($row, 18).Value2 = $r.R       # Nord
    $ws.Cells.Item($row, 19).Value2 = $r.S       # Noggrannhet
    $ws.Cells.Item($row, 20).Value2 = $r.T       # Län
    $ws.Cells.Item($row, 21).Value2 = $r.U       # Kommun
    $ws.Cells.Item($row, 22).Value2 = $r.V       # Provins
    $ws.Cells.Item($row, 23).Value2 = $r.W       # Församling
    Set-TextCell $row 25 $r.Y                    # Startdatum
    Set-TextCell $row 27 $r.AA                   # Slutdatum
    $ws.Cells.Item($row, 30).Value2 = $false     # Ej återfunnen
    $ws.Cells.Item($row, 31).Value2 = $false     # Osäker artbestämning
    Set-EmptyTextCell $row 32                    # Bestämningsmetod (blank placeholder)
    $ws.Cells.Item($row, 33).Value2 = $false     # Ospontan
    Set-EmptyTextCell $row 46                    # Bestämningsår (blank placeholder)
    $ws.Cells.Item($row, 49).Value2 = $r.AW      # Rapportör
    $ws.Cells.Item($row, 50).Value2 = $r.AX      # Observatörer
    Set-EmptyTextCell $row 51                    # Projektnamn (blank placeholder)
}
